$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 15:50"

# Country data table: row number, country name, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
  @(4, "Estados Unidos", 7681202, 1558, 4895685, 2570421, 0, 64, 215096),
  @(5, "India", 6687247, 5174, 5662490, 921128, 0, 29, 103629),
  @(6, "Brasil", 4940499, 0, 4295302, 498424, 0, 0, 146773),
  @(7, "Rusia", 1237504, 11615, 988576, 227265, 0, 188, 21663),
  @(8, "Colombia", 862158, 0, 766300, 69014, 0, 0, 26844),
  @(9, "España", 852838, 0, 0, 0, 0, 0, 32225),
  @(10, "Peru", 829999, 0, 712888, 84277, 0, 0, 32834),
  @(11, "Argentina", 809728, 0, 649017, 139243, 0, 0, 21468),
  @(12, "Mexico", 789780, 3417, 553937, 153966, 0, 180, 81877),
  @(13, "Sudafrica", 682215, 0, 615684, 49515, 0, 0, 17016),
  @(14, "Francia", 624274, 0, 98680, 493295, 0, 0, 32299),
  @(15, "Reino Unido", 515571, 0, 0, 0, 0, 0, 42369),
  @(16, "Iran", 479825, 4151, 394800, 57606, 0, 227, 27419),
  @(17, "Chile", 471746, 0, 443453, 15256, 0, 0, 13037),
  @(18, "Irak", 382949, 0, 312158, 61327, 0, 0, 9464),
  @(19, "Banglades", 371631, 1499, 284833, 81393, 0, 30, 5405),
  @(20, "Arabia Saudita", 337243, 477, 322612, 9708, 0, 25, 4923),
  @(21, "Italia", 327586, 0, 232681, 58903, 0, 0, 36002),
  @(22, "Filipinas", 326833, 2093, 273313, 47655, 0, 25, 5865),
  @(23, "Turquia", 326046, 0, 286370, 31178, 0, 0, 8498),
  @(24, "Pakistan", 315727, 467, 300616, 8588, 0, 6, 6523),
  @(25, "Indonesia", 311176, 4056, 236437, 63365, 0, 121, 11374),
  @(26, "Alemania", 305051, 394, 263700, 31727, 0, 8, 9624),
  @(27, "Israel", 274423, 2114, 208819, 63833, 0, 14, 1771),
  @(28, "Ucrania", 234584, 4348, 103401, 126663, 0, 90, 4520),
  @(29, "Canada", 168960, 0, 142334, 17122, 0, 0, 9504),
  @(30, "Paises Bajos", 144999, 4528, 0, 0, 0, 21, 6482),
  @(31, "Ecuador", 141339, 0, 120511, 9147, 0, 0, 11681),
  @(32, "Rumania", 139612, 2121, 109898, 24593, 0, 73, 5121),
  @(33, "Bolivia", 137107, 239, 98007, 30971, 0, 28, 8129),
  @(34, "Marruecos", 134695, 0, 113336, 18990, 0, 0, 2369),
  @(35, "Belgica", 132203, 1968, 19712, 102413, 0, 14, 10078),
  @(36, "Catar", 126692, 0, 123664, 2812, 0, 0, 216),
  @(37, "Panama", 115919, 0, 92423, 21066, 0, 0, 2430),
  @(38, "Republica Dominicana", 115054, 0, 90942, 21968, 0, 0, 2144),
  @(39, "Kazajistan", 108296, 60, 103367, 3204, 0, 0, 1725),
  @(40, "Kuwait", 108268, 676, 100179, 7457, 0, 4, 632),
  @(41, "Polonia", 104316, 2236, 74158, 27441, 0, 58, 2717),
  @(42, "Egipto", 103781, 0, 97398, 393, 0, 0, 5990),
  @(43, "Oman", 102648, 834, 91275, 10383, 0, 5, 990),
  @(44, "Emiratos Arabes Unidos", 100794, 1061, 90556, 9803, 0, 6, 435),
  @(45, "Suecia", 96145, 0, 0, 0, 0, 5, 5883),
  @(46, "Guatemala", 94182, 0, 82828, 8052, 0, 0, 3302),
  @(47, "Nepal", 90814, 1551, 67542, 22709, 0, 9, 563),
  @(48, "Japon", 85739, 0, 78609, 5531, 0, 0, 1599),
  @(49, "Chequia", 85566, 0, 46636, 38172, 0, 0, 758),
  @(50, "China", 85482, 12, 80635, 213, 0, 0, 4634),
  @(51, "Costa Rica", 81129, 0, 49703, 30439, 0, 0, 987),
  @(52, "Bielorrusia", 81090, 394, 75376, 4846, 0, 6, 868),
  @(53, "Honduras", 80020, 391, 29768, 47819, 0, 11, 2433),
  @(54, "Portugal", 79885, 0, 50454, 27413, 0, 0, 2018),
  @(55, "Etiopia", 79437, 0, 34016, 44191, 0, 0, 1230),
  @(56, "Venezuela", 79117, 0, 69832, 8627, 0, 0, 658),
  @(57, "Barein", 73116, 0, 67933, 4921, 0, 1, 262),
  @(58, "Nigeria", 59465, 0, 50951, 7401, 0, 0, 1113),
  @(59, "Uzbekistan", 59343, 397, 56058, 2796, 0, 4, 489),
  @(60, "Singapur", 57830, 11, 57597, 206, 0, 0, 27),
  @(61, "Moldavia", 56901, 0, 41467, 14059, 0, 0, 1375),
  @(62, "Suiza", 56632, 700, 47300, 7251, 0, 3, 2081),
  @(63, "Armenia", 53083, 406, 44932, 7161, 0, 6, 990),
  @(64, "Argelia", 52270, 0, 36672, 13830, 0, 0, 1768),
  @(65, "Austria", 49819, 923, 39790, 9207, 0, 4, 822),
  @(66, "Kirguistan", 47799, 164, 43644, 3089, 0, 0, 1066),
  @(67, "Ghana", 46829, 0, 46060, 466, 0, 0, 303),
  @(68, "Libano", 45657, 0, 20243, 25000, 0, 0, 414),
  @(69, "Paraguay", 44715, 0, 27887, 15881, 0, 0, 947),
  @(70, "Estado de Palestina", 41957, 0, 35182, 6436, 0, 0, 339),
  @(71, "Azerbaiyan", 40931, 143, 38713, 1618, 0, 2, 600),
  @(72, "Afganistan", 39486, 64, 32977, 5042, 0, 1, 1467),
  @(73, "Kenia", 39449, 0, 27035, 11679, 0, 0, 735),
  @(74, "Irlanda", 38549, 0, 23364, 13375, 0, 0, 1810),
  @(75, "Libia", 38468, 1031, 22410, 15456, 0, 6, 602),
  @(76, "Serbia", 34072, 120, 31536, 1779, 0, 1, 757),
  @(77, "Hungria", 32298, 818, 8723, 22722, 0, 20, 853),
  @(78, "Dinamarca", 30379, 322, 23655, 6061, 0, 4, 663),
  @(79, "El Salvador", 29634, 95, 24525, 4240, 0, 4, 869),
  @(80, "Bosnia y Herzegovina", 28710, 261, 22274, 5533, 0, 15, 903),
  @(81, "Australia", 27173, 24, 24892, 1386, 0, 1, 895),
  @(82, "Corea del Sur", 24239, 75, 22083, 1734, 0, 0, 422),
  @(83, "Tunez", 22230, 0, 5032, 16877, 0, 0, 321),
  @(84, "Bulgaria", 21870, 0, 15179, 5837, 0, 0, 854),
  @(85, "Camerun", 20924, 0, 19764, 740, 0, 0, 420),
  @(86, "Grecia", 20142, 0, 9989, 9736, 0, 0, 417),
  @(87, "Costa de Marfil", 19885, 0, 19490, 275, 0, 0, 120),
  @(88, "Republica de Macedonia", 18873, 0, 15487, 2626, 0, 0, 760),
  @(89, "Birmania", 18781, 0, 5548, 12789, 0, 0, 444),
  @(90, "Croacia", 18084, 287, 16192, 1588, 0, 4, 304),
  @(91, "Jordania", 17464, 0, 5292, 12062, 0, 0, 110),
  @(92, "Madagascar", 16600, 30, 15698, 668, 0, 1, 234),
  @(93, "Zambia", 15170, 81, 14313, 522, 0, 1, 335),
  @(94, "Senegal", 15141, 19, 12936, 1893, 0, 0, 312),
  @(95, "Noruega", 14669, 65, 11190, 3204, 0, 0, 275),
  @(96, "Albania", 14410, 0, 8825, 5185, 0, 0, 400),
  @(97, "Eslovaquia", 13812, 320, 5027, 8730, 0, 0, 55),
  @(98, "Sudan", 13653, 0, 6764, 6053, 0, 0, 836),
  @(99, "Malasia", 13504, 691, 10427, 2936, 0, 4, 141),
  @(100, "Montenegro", 12359, 0, 8308, 3869, 0, 0, 182),
  @(101, "Namibia", 11654, 0, 9451, 2078, 0, 0, 125),
  @(102, "Finlandia", 10929, 227, 8100, 2483, 0, 0, 346),
  @(103, "Guinea", 10800, 0, 10161, 572, 0, 0, 67),
  @(104, "Consejo Danes para los Refugiados", 10778, 0, 10239, 265, 0, 0, 274),
  @(105, "Maldivas", 10567, 0, 9427, 1106, 0, 0, 34),
  @(106, "Guayana Francesa", 10057, 0, 9710, 279, 0, 0, 68),
  @(107, "Tayikistan", 10014, 40, 8837, 1099, 0, 0, 78),
  @(108, "Mozambique", 9296, 0, 6104, 3126, 0, 0, 66),
  @(109, "Georgia", 9245, 549, 4887, 4300, 0, 4, 58),
  @(110, "Uganda", 9082, 117, 5457, 3541, 0, 2, 84),
  @(111, "Luxemburgo", 8925, 0, 7793, 1005, 0, 0, 127),
  @(112, "Haiti", 8827, 0, 6992, 1606, 0, 0, 229),
  @(113, "Gabon", 8808, 0, 8135, 619, 0, 0, 54),
  @(114, "Zimbabue", 7898, 0, 6424, 1246, 0, 0, 228),
  @(115, "Mauritania", 7523, 0, 7204, 157, 0, 0, 162),
  @(116, "Jamaica", 7012, 0, 2635, 4257, 0, 0, 120),
  @(117, "Eslovenia", 6764, 189, 4399, 2206, 0, 3, 159),
  @(118, "Cabo Verde", 6433, 0, 5524, 841, 0, 0, 68),
  @(119, "Cuba", 5845, 0, 5232, 490, 0, 0, 123),
  @(120, "Malaui", 5794, 0, 4541, 1073, 0, 0, 180),
  @(121, "Suazilandia", 5579, 0, 5141, 326, 0, 0, 112),
  @(122, "Angola", 5530, 0, 2591, 2740, 0, 0, 199),
  @(123, "Guadalupe", 5528, 0, 2199, 3272, 0, 0, 57),
  @(124, "Republica de Yibuti", 5421, 0, 5352, 8, 0, 0, 61),
  @(125, "Lituania", 5366, 81, 2546, 2721, 0, 5, 99),
  @(126, "Nicaragua", 5170, 0, 2913, 2106, 0, 0, 151),
  @(127, "Hong Kong", 5133, 8, 4875, 153, 0, 0, 105),
  @(128, "Congo", 5089, 0, 3887, 1113, 0, 0, 89),
  @(129, "Guinea Ecuatorial", 5052, 7, 4894, 75, 0, 0, 83),
  @(130, "Surinam", 4954, 0, 4755, 93, 0, 0, 106),
  @(131, "Ruanda", 4867, 0, 3226, 1612, 0, 0, 29),
  @(132, "Republica de Africa Central", 4852, 7, 1914, 2876, 0, 0, 62),
  @(133, "Trinidad yTobago", 4767, 0, 2884, 1801, 0, 0, 82),
  @(134, "Bahamas", 4452, 0, 2375, 1981, 0, 0, 96),
  @(135, "Siria", 4411, 0, 1168, 3036, 0, 0, 207),
  @(136, "Reunion", 4328, 0, 3360, 952, 0, 0, 16),
  @(137, "Aruba", 4094, 0, 3612, 451, 0, 0, 31),
  @(138, "Sri Lanka", 3979, 466, 3266, 700, 0, 0, 13),
  @(139, "Mayotte", 3892, 0, 2964, 886, 0, 0, 42),
  @(140, "Somalia", 3745, 0, 3010, 636, 0, 0, 99),
  @(141, "Estonia", 3659, 42, 2806, 786, 0, 0, 67),
  @(142, "Gambia", 3613, 19, 2233, 1263, 0, 2, 117),
  @(143, "Tailandia", 3600, 10, 3390, 151, 0, 0, 59),
  @(144, "Malta", 3374, 47, 2812, 522, 0, 1, 40),
  @(145, "Mali", 3189, 0, 2482, 576, 0, 0, 131),
  @(146, "Guyana", 3188, 0, 1972, 1126, 0, 0, 90),
  @(147, "Botsuana", 3172, 0, 710, 2446, 0, 0, 16),
  @(148, "Islandia", 3081, 101, 2324, 747, 0, 0, 10),
  @(149, "Sudan del Sur", 2726, 0, 1290, 1386, 0, 0, 50),
  @(150, "Guinea-Bisau", 2385, 0, 1728, 617, 0, 0, 40),
  @(151, "Principado de Andorra", 2370, 0, 1615, 702, 0, 0, 53),
  @(152, "Benin", 2357, 0, 1973, 343, 0, 0, 41),
  @(153, "Sierra Leona", 2269, 0, 1706, 491, 0, 0, 72),
  @(154, "Polinesia Francesa", 2228, 0, 1769, 450, 0, 0, 9),
  @(155, "Belice", 2204, 8, 1378, 796, 0, 0, 30),
  @(156, "Letonia", 2194, 68, 1322, 832, 0, 1, 40),
  @(157, "Burkina Faso", 2184, 0, 1420, 705, 0, 0, 59),
  @(158, "Uruguay", 2155, 0, 1862, 245, 0, 0, 48),
  @(159, "Yemen", 2041, 0, 1323, 126, 0, 0, 592),
  @(160, "Togo", 1864, 0, 1403, 413, 0, 0, 48),
  @(161, "Nueva Zelanda", 1858, 3, 1790, 43, 0, 0, 25),
  @(162, "Republica de Chipre", 1847, 0, 1369, 456, 0, 0, 22),
  @(163, "Lesoto", 1683, 0, 926, 718, 0, 0, 39),
  @(164, "Martinica", 1543, 0, 98, 1424, 0, 0, 21),
  @(165, "Liberia", 1354, 0, 1236, 36, 0, 0, 82),
  @(166, "Republica del Chad", 1223, 0, 1075, 62, 0, 0, 86),
  @(167, "Niger", 1200, 0, 1115, 16, 0, 0, 69),
  @(168, "Vietnam", 1098, 1, 1023, 40, 0, 0, 35),
  @(169, "Santo Tome y Principe", 913, 0, 888, 10, 0, 0, 15),
  @(170, "San Marino", 732, 0, 680, 10, 0, 0, 42),
  @(171, "Crucero", 712, 0, 651, 48, 0, 0, 13),
  @(172, "Islas Turcas y Caicos", 695, 0, 656, 33, 0, 0, 6),
  @(173, "San Martin (Parte Holandesa)", 686, 7, 592, 72, 0, 0, 22),
  @(174, "Papua Nueva Guinea", 541, 1, 527, 7, 0, 0, 7),
  @(175, "Taiwan", 521, 3, 485, 29, 0, 0, 7),
  @(176, "Burundi", 514, 0, 472, 41, 0, 0, 1),
  @(177, "Tanzania", 509, 0, 183, 305, 0, 0, 21),
  @(178, "Comoras", 487, 0, 466, 14, 0, 0, 7),
  @(179, "Islas Feroe", 475, 0, 452, 23, 0, 0, 0),
  @(180, "Curazao", 462, 0, 242, 219, 0, 0, 1),
  @(181, "Gibraltar", 437, 5, 368, 69, 0, 0, 0),
  @(182, "San Martin (Parte Francesa)", 412, 0, 309, 95, 0, 0, 8),
  @(183, "Eritrea", 398, 0, 358, 40, 0, 0, 0),
  @(184, "Mauricio", 387, 0, 357, 20, 0, 0, 10),
  @(185, "Isla de Man", 344, 0, 316, 4, 0, 0, 24),
  @(186, "Mongolia", 315, 1, 307, 8, 0, 0, 0),
  @(187, "Butan", 299, 1, 248, 51, 0, 0, 0),
  @(188, "Camboya", 280, 0, 276, 4, 0, 0, 0),
  @(189, "Monaco", 223, 0, 193, 28, 0, 0, 2),
  @(190, "Islas Caimanes", 213, 0, 210, 2, 0, 0, 1),
  @(191, "Barbados", 200, 0, 182, 11, 0, 0, 7),
  @(192, "Bermudas", 181, 0, 170, 2, 0, 0, 9),
  @(193, "Seychelles", 146, 0, 143, 3, 0, 0, 0),
  @(194, "Brunei", 146, 0, 143, 0, 0, 0, 3),
  @(195, "Bonaire, San Eustaquio y Saba", 141, 17, 67, 72, 0, 0, 2),
  @(196, "Liechtenstein", 130, 3, 116, 13, 0, 0, 1),
  @(197, "Antigua y Barbuda", 107, 0, 96, 8, 0, 0, 3),
  @(198, "Islas Virgenes Britanicas", 71, 0, 66, 4, 0, 0, 1),
  @(199, "San Vicente y las Granadinas", 64, 0, 64, 0, 0, 0, 0),
  @(200, "San Bartolome", 62, 0, 37, 25, 0, 0, 0),
  @(201, "Macao", 46, 0, 46, 0, 0, 0, 0),
  @(202, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @(203, "Guam", 32, 0, 0, 31, 0, 0, 1),
  @(204, "Fiyi", 32, 0, 28, 2, 0, 0, 2),
  @(205, "Dominica", 31, 0, 24, 7, 0, 0, 0),
  @(206, "Timor Oriental", 28, 0, 28, 0, 0, 0, 0),
  @(207, "Nueva Caledonia", 27, 0, 27, 0, 0, 0, 0),
  @(208, "Santa Lucia", 27, 0, 27, 0, 0, 0, 0),
  @(209, "Granada", 24, 0, 24, 0, 0, 0, 0),
  @(210, "Laos", 23, 0, 22, 1, 0, 0, 0),
  @(211, "San Cristobal y Nieves", 19, 0, 17, 2, 0, 0, 0),
  @(212, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @(213, "San Pedro y Miquelon", 16, 0, 6, 10, 0, 0, 0),
  @(214, "Groenlandia", 14, 0, 14, 0, 0, 0, 0),
  @(215, "Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
  @(216, "Montserrat", 13, 0, 12, 0, 0, 0, 1),
  @(217, "Santa Sede", 12, 0, 12, 0, 0, 0, 0),
  @(218, "Sahara Occidental", 10, 0, 8, 1, 0, 0, 1),
  @(219, "Anguila", 3, 0, 3, 0, 0, 0, 0),
  @(220, "Islas Salomon", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
